$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the temporary sample data that was used for automatic deployment
# testing, leaving the date number-format styling intact on column A.
$ws.Range("A2:F7").ClearContents()

# Reproduce the recorded selection: the user highlighted the block A2:G8
# (selecting from G8 back up to A2) before deleting it.
$ws.Range("A2:G8").Select()
